$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "70.705.13"
$ws.Cells.Item(2, 5).Value = "  +1.03%  "
$ws.Cells.Item(3, 4).Value = "3.527.75"
$ws.Cells.Item(3, 5).Value = "  +0.14%  "
$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.07%  "
$ws.Cells.Item(5, 4).Value = "'606.60"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.07%  "
$ws.Cells.Item(6, 4).Value = "'175.04"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +1.98%  "
$ws.Cells.Item(7, 4).Value = "'0.613"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -0.78%  "
$ws.Cells.Item(8, 4).Value = "3.523.54"
$ws.Cells.Item(8, 5).Value = "  +0.15%  "
$ws.Cells.Item(9, 5).Value = "  +0.02%  "
$ws.Cells.Item(10, 4).Value = "'0.195"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -2.57%  "
$ws.Cells.Item(11, 4).Value = "'7.22"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +7.56%  "
$ws.Cells.Item(12, 4).Value = "'0.587"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +0.50%  "
$ws.Cells.Item(13, 4).Value = "'46.47"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -1.93%  "
$ws.Cells.Item(14, 4).Value = "'0.0000277"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -0.83%  "
$ws.Cells.Item(15, 4).Value = "4.100.59"
$ws.Cells.Item(15, 5).Value = "  +0.27%  "
$ws.Cells.Item(16, 4).Value = "'8.33"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -0.66%  "
$ws.Cells.Item(17, 4).Value = "'612.08"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -1.74%  "
$ws.Cells.Item(18, 2).Value = "WrappedBTC"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(18, 4).Value = "70.723.68"
$ws.Cells.Item(18, 5).Value = "  +0.95%  "
$ws.Cells.Item(19, 2).Value = "WrappedEther"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(19, 4).Value = "3.520.87"
$ws.Cells.Item(19, 5).Value = "  -0.12%  "
$ws.Cells.Item(20, 5).Value = "  +0.87%  "
$ws.Cells.Item(21, 4).Value = "'17.58"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +1.43%  "
$ws.Cells.Item(22, 4).Value = "'0.882"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -0.45%  "
$ws.Cells.Item(23, 5).Value = "  -9.13%  "
$ws.Cells.Item(24, 4).Value = "'98.99"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +2.96%  "
$ws.Cells.Item(25, 4).Value = "'15.66"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -1.11%  "
$ws.Cells.Item(26, 5).Value = "  -3.44%  "
$ws.Cells.Item(28, 4).Value = "'2.58"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -1.41%  "
$ws.Cells.Item(29, 4).Value = "'34.07"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +2.54%  "
$ws.Cells.Item(30, 4).Value = "'9.06"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -2.52%  "
$ws.Cells.Item(31, 5).Value = "  -3.31%  "
$ws.Cells.Item(32, 4).Value = "'8.09"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -4.56%  "
$ws.Cells.Item(33, 4).Value = "'640.43"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +12.65%  "
$ws.Cells.Item(34, 5).Value = "  -4.39%  "
$ws.Cells.Item(35, 4).Value = "'6.85"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -2.17%  "
$ws.Cells.Item(36, 4).Value = "'3.59"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +1.75%  "
$ws.Cells.Item(37, 4).Value = "'0.0998"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -1.93%  "
$ws.Cells.Item(38, 4).Value = "'10.79"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -0.05%  "
$ws.Cells.Item(39, 4).Value = "'0.0478"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +5.50%  "
$ws.Cells.Item(40, 4).Value = "'56.85"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -0.23%  "
$ws.Cells.Item(41, 4).Value = "'0.998"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -0.20%  "
$ws.Cells.Item(42, 4).Value = "'0.142"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +0.86%  "
$ws.Cells.Item(43, 2).Value = "PEPE"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(43, 4).Value = "0.0₃0744"
$ws.Cells.Item(43, 5).Value = "  +5.50%  "
$ws.Cells.Item(44, 2).Value = "Maker"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(44, 4).Value = "3.378.91"
$ws.Cells.Item(44, 5).Value = "  +1.12%  "
$ws.Cells.Item(45, 4).Value = "'0.311"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -5.25%  "
$ws.Cells.Item(46, 4).Value = "'2.91"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -2.40%  "
$ws.Cells.Item(47, 4).Value = "'32.25"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -2.81%  "
$ws.Cells.Item(48, 4).Value = "'2.57"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -3.02%  "
$ws.Cells.Item(49, 5).Value = "  +0.84%  "
$ws.Cells.Item(50, 4).Value = "'133.42"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -2.08%  "
$ws.Cells.Item(51, 5).Value = "  -0.02%  "
